$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row for Perejil (Vega Central Mapocho de Santiago)
# is inserted at row 183, pushing the existing rows 183:220 down to 184:221.
$ws.Rows.Item(183).Insert()

$ws.Cells.Item(183, 1).Value = 9
$ws.Cells.Item(183, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(183, 3).Value = "Metropolitana"
$ws.Cells.Item(183, 4).Value = 44476
$ws.Cells.Item(183, 5).Value = 13
$ws.Cells.Item(183, 6).Value = 100112044
$ws.Cells.Item(183, 7).Value = "Perejil"
$ws.Cells.Item(183, 8).Value = "Sin especificar"
$ws.Cells.Item(183, 9).Value = "Primera"
$ws.Cells.Item(183, 10).Value = 97
$ws.Cells.Item(183, 11).Value = 8000
$ws.Cells.Item(183, 12).Value = 10000
$ws.Cells.Item(183, 13).Value = 8990
$ws.Cells.Item(183, 14).Value = "$/docena de atados"
$ws.Cells.Item(183, 15).Value = "Región Metropolitana"
$ws.Cells.Item(183, 16).Value = 2997
$ws.Cells.Item(183, 17).Value = 3
$ws.Cells.Item(183, 18).Value = "Hortaliza"
